$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("Area") values were recomputed (e.g. pixel counts converted to
# calibrated area units) for rows 2-13.
$ws.Range("B2").Value = 175.774
$ws.Range("B3").Value = 939.46400000000006
$ws.Range("B4").Value = 590.01099999999997
$ws.Range("B5").Value = 207.483
$ws.Range("B6").Value = 482.54899999999998
$ws.Range("B7").Value = 231.69800000000001
$ws.Range("B8").Value = 912.803
$ws.Range("B9").Value = 260.55700000000002
$ws.Range("B10").Value = 200.95099999999999
$ws.Range("B11").Value = 273.98500000000001
$ws.Range("B12").Value = 248.756
$ws.Range("B13").Value = 241.99

# The sheet's lingering "O16" selection (outside the used range A1:C13) is
# cleared back to the default top-left cell.
[void]$ws.Range("A1").Select()
